# Generate Report for Handoff
# Adds two new files (4f2d23ac..., c7fcc05d...) to the localization status
# report and refreshes the existing 65688d7c... row with its new handoff
# timestamps, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ovTbl = $ov.ListObjects.Item(1)

# Row 3 (existing row) now describes the newly-added 4f2d23ac file.
$ov.Range("A3").Value = "4f2d23ac-c938-4c84-9351-0034a1c0dd8a.md"
$ov.Range("C3").Value = ".md"
$ov.Range("D3").Value = ""
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-09-02 06:19:15"
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5db86f7c6342233b7be3c1a0ffb61e0a61e90f59/e2e/f9c17345-6c55-4529-be4c-7456e6e48e97.md", $null, $null, "e2e\f9c17345-6c55-4529-be4c-7456e6e48e97.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5ccf3878fe38fa813f6d699e240a24998449a6ab/e2e/4f2d23ac-c938-4c84-9351-0034a1c0dd8a.md", $null, $null, "e2e\4f2d23ac-c938-4c84-9351-0034a1c0dd8a.md") | Out-Null

# Row 4 (new) - the 65688d7c file that used to sit in row 3.
$ovRow4 = $ovTbl.ListRows.Add()
$ov.Range("A4").Value = "65688d7c-20fa-4af5-9068-cbe2e4639b2c.md"
$ov.Range("C4").Value = ".md"
$ov.Range("D4").Value = ""
$ov.Range("E4").Value = "Ready for handoff"
$ov.Range("F4").Value = "Ready for handoff"
$ov.Range("G4").Value = "2016-09-02 06:19:15"
$ov.Hyperlinks.Add($ov.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a38e86403e12fe00aa0d5d77cb9b8c1b7755d05/e2e/65688d7c-20fa-4af5-9068-cbe2e4639b2c.md", $null, $null, "e2e\65688d7c-20fa-4af5-9068-cbe2e4639b2c.md") | Out-Null

# Row 5 (new) - the brand new c7fcc05d file.
$ovRow5 = $ovTbl.ListRows.Add()
$ov.Range("A5").Value = "c7fcc05d-1f59-4621-99b1-649c6a5a2de4.md"
$ov.Range("C5").Value = ".md"
$ov.Range("D5").Value = ""
$ov.Range("E5").Value = "Ready for handoff"
$ov.Range("F5").Value = "Ready for handoff"
$ov.Range("G5").Value = "2016-09-02 06:19:15"
$ov.Hyperlinks.Add($ov.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87820044c8ae2c779c6f1731be196907c45bbbfa/e2e/c7fcc05d-1f59-4621-99b1-649c6a5a2de4.md", $null, $null, "e2e\c7fcc05d-1f59-4621-99b1-649c6a5a2de4.md") | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zhTbl = $zh.ListObjects.Item(1)

# Row 3 (existing row) -> 4f2d23ac
$zh.Range("A3").Value = "4f2d23ac-c938-4c84-9351-0034a1c0dd8a.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "e2e"
$zh.Range("E3").Value = ""
$zh.Range("F3").Value = "False"
$zh.Range("G3").Value = "4f2d23ac-c938-4c84-9351-0034a1c0dd8a.e62e23a0058dffad3f2e0c236f4034754d5c2742.zh-cn.xlf"
$zh.Range("H3").Value = "2016-09-02 06:19:10"
$zh.Range("I3").Value = ""
$zh.Range("J3").Value = ""
$zh.Range("K3").Value = "0001-01-01 00:00:00"
$zh.Range("L3").Value = ""
$zh.Range("M3").Value = "True"
$zh.Range("N3").Value = ""
$zh.Range("O3").Value = "False"
$zh.Range("P3").Value = ""
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5db86f7c6342233b7be3c1a0ffb61e0a61e90f59/e2e/f9c17345-6c55-4529-be4c-7456e6e48e97.md", $null, $null, "f9c17345-6c55-4529-be4c-7456e6e48e97.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5ccf3878fe38fa813f6d699e240a24998449a6ab/e2e/4f2d23ac-c938-4c84-9351-0034a1c0dd8a.md", $null, $null, "4f2d23ac-c938-4c84-9351-0034a1c0dd8a.md") | Out-Null

# Row 4 (new) -> 65688d7c
$zhRow4 = $zhTbl.ListRows.Add()
$zh.Range("A4").Value = "65688d7c-20fa-4af5-9068-cbe2e4639b2c.md"
$zh.Range("B4").Value = ".md"
$zh.Range("C4").Value = "Ready for handoff"
$zh.Range("D4").Value = "e2e"
$zh.Range("E4").Value = ""
$zh.Range("F4").Value = "False"
$zh.Range("G4").Value = "65688d7c-20fa-4af5-9068-cbe2e4639b2c.93e43df5a0e11ea6cb4405509607e0678164e1de.zh-cn.xlf"
$zh.Range("H4").Value = "2016-09-02 06:19:10"
$zh.Range("I4").Value = ""
$zh.Range("J4").Value = ""
$zh.Range("K4").Value = "0001-01-01 00:00:00"
$zh.Range("L4").Value = ""
$zh.Range("M4").Value = "True"
$zh.Range("N4").Value = ""
$zh.Range("O4").Value = "False"
$zh.Range("P4").Value = ""
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a38e86403e12fe00aa0d5d77cb9b8c1b7755d05/e2e/65688d7c-20fa-4af5-9068-cbe2e4639b2c.md", $null, $null, "65688d7c-20fa-4af5-9068-cbe2e4639b2c.md") | Out-Null

# Row 5 (new) -> c7fcc05d
$zhRow5 = $zhTbl.ListRows.Add()
$zh.Range("A5").Value = "c7fcc05d-1f59-4621-99b1-649c6a5a2de4.md"
$zh.Range("B5").Value = ".md"
$zh.Range("C5").Value = "Ready for handoff"
$zh.Range("D5").Value = "e2e"
$zh.Range("E5").Value = ""
$zh.Range("F5").Value = "False"
$zh.Range("G5").Value = "c7fcc05d-1f59-4621-99b1-649c6a5a2de4.30637ce979c95581bbc51cd880944de99fd214da.zh-cn.xlf"
$zh.Range("H5").Value = "2016-09-02 06:19:10"
$zh.Range("I5").Value = ""
$zh.Range("J5").Value = ""
$zh.Range("K5").Value = "0001-01-01 00:00:00"
$zh.Range("L5").Value = ""
$zh.Range("M5").Value = "True"
$zh.Range("N5").Value = ""
$zh.Range("O5").Value = "False"
$zh.Range("P5").Value = ""
$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87820044c8ae2c779c6f1731be196907c45bbbfa/e2e/c7fcc05d-1f59-4621-99b1-649c6a5a2de4.md", $null, $null, "c7fcc05d-1f59-4621-99b1-649c6a5a2de4.md") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$deTbl = $de.ListObjects.Item(1)

# Row 3 (existing row) -> 4f2d23ac
$de.Range("A3").Value = "4f2d23ac-c938-4c84-9351-0034a1c0dd8a.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "e2e"
$de.Range("E3").Value = ""
$de.Range("F3").Value = "False"
$de.Range("G3").Value = "4f2d23ac-c938-4c84-9351-0034a1c0dd8a.e62e23a0058dffad3f2e0c236f4034754d5c2742.de-de.xlf"
$de.Range("H3").Value = "2016-09-02 06:19:15"
$de.Range("I3").Value = ""
$de.Range("J3").Value = ""
$de.Range("K3").Value = "0001-01-01 00:00:00"
$de.Range("L3").Value = ""
$de.Range("M3").Value = "True"
$de.Range("N3").Value = ""
$de.Range("O3").Value = "False"
$de.Range("P3").Value = ""
$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5db86f7c6342233b7be3c1a0ffb61e0a61e90f59/e2e/f9c17345-6c55-4529-be4c-7456e6e48e97.md", $null, $null, "f9c17345-6c55-4529-be4c-7456e6e48e97.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5ccf3878fe38fa813f6d699e240a24998449a6ab/e2e/4f2d23ac-c938-4c84-9351-0034a1c0dd8a.md", $null, $null, "4f2d23ac-c938-4c84-9351-0034a1c0dd8a.md") | Out-Null

# Row 4 (new) -> 65688d7c
$deRow4 = $deTbl.ListRows.Add()
$de.Range("A4").Value = "65688d7c-20fa-4af5-9068-cbe2e4639b2c.md"
$de.Range("B4").Value = ".md"
$de.Range("C4").Value = "Ready for handoff"
$de.Range("D4").Value = "e2e"
$de.Range("E4").Value = ""
$de.Range("F4").Value = "False"
$de.Range("G4").Value = "65688d7c-20fa-4af5-9068-cbe2e4639b2c.93e43df5a0e11ea6cb4405509607e0678164e1de.de-de.xlf"
$de.Range("H4").Value = "2016-09-02 06:19:15"
$de.Range("I4").Value = ""
$de.Range("J4").Value = ""
$de.Range("K4").Value = "0001-01-01 00:00:00"
$de.Range("L4").Value = ""
$de.Range("M4").Value = "True"
$de.Range("N4").Value = ""
$de.Range("O4").Value = "False"
$de.Range("P4").Value = ""
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a38e86403e12fe00aa0d5d77cb9b8c1b7755d05/e2e/65688d7c-20fa-4af5-9068-cbe2e4639b2c.md", $null, $null, "65688d7c-20fa-4af5-9068-cbe2e4639b2c.md") | Out-Null

# Row 5 (new) -> c7fcc05d
$deRow5 = $deTbl.ListRows.Add()
$de.Range("A5").Value = "c7fcc05d-1f59-4621-99b1-649c6a5a2de4.md"
$de.Range("B5").Value = ".md"
$de.Range("C5").Value = "Ready for handoff"
$de.Range("D5").Value = "e2e"
$de.Range("E5").Value = ""
$de.Range("F5").Value = "False"
$de.Range("G5").Value = "c7fcc05d-1f59-4621-99b1-649c6a5a2de4.30637ce979c95581bbc51cd880944de99fd214da.de-de.xlf"
$de.Range("H5").Value = "2016-09-02 06:19:15"
$de.Range("I5").Value = ""
$de.Range("J5").Value = ""
$de.Range("K5").Value = "0001-01-01 00:00:00"
$de.Range("L5").Value = ""
$de.Range("M5").Value = "True"
$de.Range("N5").Value = ""
$de.Range("O5").Value = "False"
$de.Range("P5").Value = ""
$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87820044c8ae2c779c6f1731be196907c45bbbfa/e2e/c7fcc05d-1f59-4621-99b1-649c6a5a2de4.md", $null, $null, "c7fcc05d-1f59-4621-99b1-649c6a5a2de4.md") | Out-Null
